$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The label for 01-4.png (row 5, column C) was edited: the trailing
# "화난 앉아있는 " (with a trailing space) became "화난엄마 앉아있는"
# (an "엄마" inserted, trailing space dropped).
$ws.Range("C5").Value = "2인 거실 책읽기 책 어질러진 화난 5살남자아이검은색티셔츠금색짧은머리금발화난표정화난 30살검은색긴머리흰색티셔츠화난표정화난엄마 앉아있는"

# The active selection moved from C7 to C5.
$ws.Range("C5").Select()
